# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted as row 31 (pushing the existing
# rows 31..107 down to 32..108). The new row carries a later date
# (2023-01-26, serial 44952) together with its own volume/price figures,
# while all of the descriptive (constant) columns keep the same values used
# throughout the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 31; this shifts rows 31-107
# down to 32-108 without touching their contents.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new record.
$ws.Range("A31").Value = 2
$ws.Range("B31").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = "2023-01-26"
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 100112030
$ws.Range("G31").Value = "Poroto granado"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 600
$ws.Range("K31").Value = 25000
$ws.Range("L31").Value = 27000
$ws.Range("M31").Value = 26000
$ws.Range("N31").Value = "`$/malla 25 kilos"
$ws.Range("O31").Value = "Provincia de Limarí"
$ws.Range("P31").Value = 1040
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"
